$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3646.5
$ws.Range("J18").Value = 7164.3335
$ws.Range("L18").Value = 7164.3335
$ws.Range("N18").Value = -7732.3335
$ws.Range("H33").Value = 270.91666
$ws.Range("I33").Value = 252.44444
$ws.Range("K33").Value = 252.44444
$ws.Range("M33").Value = -23.44443999999999
$ws.Range("H116").Value = 8160.615
$ws.Range("I116").Value = 8427.5
$ws.Range("K116").Value = 8427.5
$ws.Range("M116").Value = -4985.5
$ws.Range("H137").Value = 2678.6875
$ws.Range("I137").Value = 2541.111
$ws.Range("J137").Value = 2855.5715
$ws.Range("K137").Value = 7623.333
$ws.Range("L137").Value = 8566.7145
$ws.Range("M137").Value = -5073.333
$ws.Range("N137").Value = -13666.7145
$ws.Range("H138").Value = 1485715.4
$ws.Range("I138").Value = 5898.2085
$ws.Range("J138").Value = 2182100
$ws.Range("K138").Value = 17694.6255
$ws.Range("L138").Value = 6546300
$ws.Range("M138").Value = -12554.6255
$ws.Range("N138").Value = -6556580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7620.6733
$ws.Range("I32").Value = 6160.5
$ws.Range("K32").Value = 6160.5
$ws.Range("M32").Value = -5873.5
$ws.Range("H37").Value = 23011.25
$ws.Range("J37").Value = 39997.5
$ws.Range("L37").Value = 39997.5
$ws.Range("N37").Value = -40543.5
$ws.Range("H45").Value = 6549.5713
$ws.Range("I45").Value = 8976.691999999999
$ws.Range("K45").Value = 8976.691999999999
$ws.Range("M45").Value = -8599.691999999999
$ws.Range("H132").Value = 2439.1538
$ws.Range("I132").Value = 2184.6667
$ws.Range("J132").Value = 3011.75
$ws.Range("K132").Value = 6554.000100000001
$ws.Range("L132").Value = 9035.25
$ws.Range("M132").Value = -4024.000100000001
$ws.Range("N132").Value = -14095.25
$ws.Range("H133").Value = 60808.375
$ws.Range("J133").Value = 60808.375
$ws.Range("L133").Value = 60808.375
$ws.Range("N133").Value = -65868.375
$ws.Range("H135").Value = 230357
$ws.Range("J135").Value = 230357
$ws.Range("L135").Value = 230357
$ws.Range("N135").Value = -240497
$ws.Range("H139").Value = 243749.75
$ws.Range("J139").Value = 243749.75
$ws.Range("L139").Value = 243749.75
$ws.Range("N139").Value = -254029.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4931.875
$ws.Range("J94").Value = 6665.3335
$ws.Range("L94").Value = 6665.3335
$ws.Range("N94").Value = -7567.3335
$ws.Range("H105").Value = 6848.032
$ws.Range("I105").Value = 8175.476
$ws.Range("K105").Value = 8175.476
$ws.Range("M105").Value = -6428.476
$ws.Range("H107").Value = 2944.2964
$ws.Range("I107").Value = 1936.2046
$ws.Range("K107").Value = 1936.2046
$ws.Range("M107").Value = -16.20460000000003
$ws.Range("H134").Value = 2428.739
$ws.Range("I134").Value = 2374.3333
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 7122.999899999999
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4587.999899999999
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 701.25
$ws.Range("I19").Value = 701.25
$ws.Range("K19").Value = 701.25
$ws.Range("M19").Value = -531.25
$ws.Range("H24").Value = 701.25
$ws.Range("I24").Value = 701.25
$ws.Range("K24").Value = 701.25
$ws.Range("M24").Value = -531.25
$ws.Range("H58").Value = 3362.2778
$ws.Range("I58").Value = 3401.3125
$ws.Range("J58").Value = 3050
$ws.Range("K58").Value = 3401.3125
$ws.Range("L58").Value = 3050
$ws.Range("M58").Value = -3198.3125
$ws.Range("N58").Value = -3456
$ws.Range("H62").Value = 6642.125
$ws.Range("J62").Value = 7813.2856
$ws.Range("L62").Value = 7813.2856
$ws.Range("N62").Value = -9061.285599999999
$ws.Range("H65").Value = 6642.125
$ws.Range("J65").Value = 7813.2856
$ws.Range("L65").Value = 39066.428
$ws.Range("N65").Value = -45306.428
$ws.Range("H99").Value = 2156.95
$ws.Range("I99").Value = 1946
$ws.Range("K99").Value = 1946
$ws.Range("M99").Value = -448
$ws.Range("H126").Value = 2156.95
$ws.Range("I126").Value = 1946
$ws.Range("K126").Value = 5838
$ws.Range("M126").Value = -3368
$ws.Range("H134").Value = 15158
$ws.Range("I134").Value = 15665.182
$ws.Range("K134").Value = 46995.546
$ws.Range("M134").Value = -44460.546
$ws.Range("H136").Value = 3362.2778
$ws.Range("I136").Value = 3401.3125
$ws.Range("J136").Value = 3050
$ws.Range("K136").Value = 10203.9375
$ws.Range("L136").Value = 9150
$ws.Range("M136").Value = -7653.9375
$ws.Range("N136").Value = -14250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10735.25
$ws.Range("J11").Value = 400
$ws.Range("L11").Value = 1200
$ws.Range("N11").Value = -1480
$ws.Range("H107").Value = 1139.6
$ws.Range("I107").Value = 1435.5
$ws.Range("J107").Value = 1000.35297
$ws.Range("K107").Value = 4306.5
$ws.Range("L107").Value = 3001.05891
$ws.Range("M107").Value = -2386.5
$ws.Range("N107").Value = -6841.05891
$ws.Range("H139").Value = 6943.1665
$ws.Range("I139").Value = 6943.1665
$ws.Range("K139").Value = 20829.4995
$ws.Range("M139").Value = -15689.4995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 9242.272000000001
$ws.Range("I2").Value = 14337.286
$ws.Range("J2").Value = 326
$ws.Range("K2").Value = 14337.286
$ws.Range("L2").Value = 326
$ws.Range("M2").Value = -14224.286
$ws.Range("N2").Value = -552
$ws.Range("H11").Value = 2576012
$ws.Range("J11").Value = 66010.22
$ws.Range("L11").Value = 66010.22
$ws.Range("N11").Value = -66288.22
$ws.Range("H107").Value = 1619.2632
$ws.Range("I107").Value = 553.7143
$ws.Range("J107").Value = 2240.8333
$ws.Range("K107").Value = 553.7143
$ws.Range("L107").Value = 2240.8333
$ws.Range("M107").Value = 1366.2857
$ws.Range("N107").Value = -6080.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H55").Value = 2189.4
$ws.Range("I55").Value = 1732.8334
$ws.Range("K55").Value = 1732.8334
$ws.Range("M55").Value = -1559.8334
$ws.Range("H122").Value = 5953.84
$ws.Range("I122").Value = 5849.8096
$ws.Range("K122").Value = 17549.4288
$ws.Range("M122").Value = -15099.4288
$ws.Range("H132").Value = 3541.2083
$ws.Range("I132").Value = 2819.4
$ws.Range("K132").Value = 8458.200000000001
$ws.Range("M132").Value = -5928.200000000001
$ws.Range("H140").Value = 54171.6
$ws.Range("J140").Value = 54171.6
$ws.Range("L140").Value = 54171.6
$ws.Range("N140").Value = -64531.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 33600000
$ws.Range("J5").Value = 50150000
$ws.Range("L5").Value = 50150000
$ws.Range("N5").Value = -50150224
$ws.Range("H62").Value = 10872.682
$ws.Range("I62").Value = 7457
$ws.Range("K62").Value = 7457
$ws.Range("M62").Value = -6833
$ws.Range("H65").Value = 10872.682
$ws.Range("I65").Value = 7457
$ws.Range("K65").Value = 37285
$ws.Range("M65").Value = -34165
$ws.Range("H113").Value = 859.4516
$ws.Range("I113").Value = 783.9
$ws.Range("J113").Value = 996.8182
$ws.Range("K113").Value = 2351.7
$ws.Range("L113").Value = 2990.4546
$ws.Range("M113").Value = -181.6999999999998
$ws.Range("N113").Value = -7330.4546
$ws.Range("H122").Value = 7490.4546
$ws.Range("I122").Value = 2199
$ws.Range("K122").Value = 6597
$ws.Range("M122").Value = -4147
$ws.Range("H132").Value = 4780.143
$ws.Range("I132").Value = 3568.5
$ws.Range("J132").Value = 7203.4287
$ws.Range("K132").Value = 10705.5
$ws.Range("L132").Value = 21610.2861
$ws.Range("M132").Value = -8175.5
$ws.Range("N132").Value = -26670.2861
